# Updated symbol list on Wed Feb  8 19:32:59 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking table. Values are written as literal text (leading "'"
# forces text storage, matching the workbook's existing inlineStr cells)
# so cells keep their original General/quote-prefixed style rather than
# being reinterpreted as numbers or percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.18"
$ws.Range("E2").Value = "'-0.61%"
$ws.Range("D3").Value = "'44.49"
$ws.Range("E3").Value = "'0.64%"
$ws.Range("D4").Value = "'5.206"
$ws.Range("E4").Value = "'-5.67%"
$ws.Range("D5").Value = "'0.08377"
$ws.Range("E5").Value = "'3.45%"
$ws.Range("D6").Value = "'1.939"
$ws.Range("E6").Value = "'-5.95%"
$ws.Range("D7").Value = "'0.9727"
$ws.Range("E7").Value = "'0.13%"
$ws.Range("D8").Value = "'2.505"
$ws.Range("E8").Value = "'-5.60%"
$ws.Range("D9").Value = "'0.1154"
$ws.Range("E9").Value = "'3.59%"
$ws.Range("D10").Value = "'0.1907"
$ws.Range("E10").Value = "'1.24%"
$ws.Range("D11").Value = "'0.09659"
$ws.Range("E11").Value = "'-3.15%"
$ws.Range("E12").Value = "'-3.35%"
$ws.Range("E13").Value = "'0.28%"
$ws.Range("D14").Value = "'0.001291"
$ws.Range("E14").Value = "'1.61%"
$ws.Range("D15").Value = "'0.005823"
$ws.Range("E15").Value = "'-3.07%"
$ws.Range("D16").Value = "'3.402"
$ws.Range("E16").Value = "'1.83%"
$ws.Range("D17").Value = "'4.447"
$ws.Range("E17").Value = "'0.31%"
$ws.Range("E18").Value = "'1.74%"
$ws.Range("D19").Value = "'8.675"
$ws.Range("E19").Value = "'-14.96%"
$ws.Range("D20").Value = "'0.1363"
$ws.Range("E20").Value = "'-2.06%"
$ws.Range("E21").Value = "'0.29%"
$ws.Range("D22").Value = "'0.04150"
$ws.Range("E22").Value = "'1.38%"
$ws.Range("E23").Value = "'-5.54%"
$ws.Range("D24").Value = "'0.004433"
$ws.Range("E24").Value = "'1.07%"
$ws.Range("D25").Value = "'0.0001304"
$ws.Range("E25").Value = "'1.96%"
$ws.Range("D26").Value = "'0.0002986"
$ws.Range("E26").Value = "'-20.16%"
$ws.Range("D38").Value = "'0.02740"
$ws.Range("E38").Value = "'2.72%"
$ws.Range("E39").Value = "'0.11%"
$ws.Range("D40").Value = "'0.007874"
$ws.Range("E40").Value = "'3.77%"
$ws.Range("D41").Value = "'0.1412"
$ws.Range("E41").Value = "'0.10%"
$ws.Range("D42").Value = "'0.007309"
$ws.Range("E42").Value = "'-11.33%"
$ws.Range("D43").Value = "'0.002046"
$ws.Range("E43").Value = "'4.56%"
$ws.Range("D44").Value = "'0.007926"
$ws.Range("D45").Value = "'0.3502"
$ws.Range("E46").Value = "'-3.56%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.34%"
$ws.Range("D48").Value = "'0.003495"
$ws.Range("E48").Value = "'-3.79%"
$ws.Range("D49").Value = "'0.003538"
$ws.Range("E49").Value = "'40.40%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.34%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.34%"
